# feat: add 2022-Q1 data
# - Inserts a new "2022-Q1" worksheet (fund holdings detail) positioned
#   right after "2021-Q4" and right before "总计".
# - Rebuilds the "总计" (totals) summary worksheet with a new first data
#   row for 2022-Q1 (16 holdings, 8.72 亿元), shifting the previously
#   existing rows down by one.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ---------------------------------------------------------------------
# 1. Drop "总计" so it can be re-created *after* the new 2022-Q1 sheet --
#    this keeps sheet/tab order = [...,"2021-Q4","2022-Q1","总计"].
#    (Its rows are rebuilt from scratch below, with 2022-Q1 prepended.)
# ---------------------------------------------------------------------
$wb.Worksheets.Item("总计").Delete()

# ---------------------------------------------------------------------
# 2. Create "2022-Q1" right after "2021-Q4".
# ---------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("2021-Q4")
$q1Sheet = $wb.Worksheets.Add($null, $afterSheet)
$q1Sheet.Name = "2022-Q1"

$q1Headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 0; $col -lt $q1Headers.Length; $col++) {
    $cell = $q1Sheet.Cells.Item(1, $col + 2)
    $cell.Value = $q1Headers[$col]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$q1Data = @(
    @("162703", "广发小盘成长混合(LOF)A", "98.71", "94.88", "3.52", "3.4746", 10),
    @("003745", "广发多元新兴股票", "37.08", "90.97", "5.40", "2.0023", 7),
    @("005299", "万家成长优选灵活配置混合A", "24.43", "93.91", "2.69", "0.6572", 8),
    @("010694", "万家内需增长一年持有期混合", "17.21", "94.85", "3.58", "0.6161", 7),
    @("001239", "长盛国企改革主题灵活配置混合", "5.02", "87.24", "8.70", "0.4367", 2),
    @("010611", "万家战略发展产业混合A", "12.80", "93.73", "2.91", "0.3725", 10),
    @("005300", "万家成长优选灵活配置混合C", "10.93", "93.91", "2.69", "0.2940", 8),
    @("510081", "长盛动态精选混合", "3.15", "60.76", "6.34", "0.1997", 1),
    @("009132", "广发小盘成长混合(LOF)C", "5.31", "94.88", "3.52", "0.1869", 10),
    @("006132", "万家智造优势混合A", "4.82", "93.70", "3.86", "0.1861", 4),
    @("010612", "万家战略发展产业混合C", "5.06", "93.73", "2.91", "0.1472", 10),
    @("970043", "东吴裕盈一年持有期灵活配置混合A", "2.06", "57.66", "2.92", "0.0602", 6),
    @("000354", "长盛城镇化主题混合", "0.43", "78.41", "8.51", "0.0366", 1),
    @("006133", "万家智造优势混合C", "0.52", "93.70", "3.86", "0.0201", 4),
    @("970045", "东吴裕盈一年持有期灵活配置混合C", "0.56", "57.66", "2.92", "0.0164", 6),
    @("970044", "东吴裕盈一年持有期灵活配置混合B", "0.31", "57.66", "2.92", "0.0091", 6)
)

for ($i = 0; $i -lt $q1Data.Length; $i++) {
    $r = $i + 2
    $row = $q1Data[$i]

    $aCell = $q1Sheet.Cells.Item($r, 1)
    $aCell.Value = $i
    $aCell.Font.Bold = $true
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160
    $aCell.Borders.LineStyle = 1

    for ($col = 0; $col -lt 6; $col++) {
        $cell = $q1Sheet.Cells.Item($r, $col + 2)
        $cell.NumberFormat = "@"
        $cell.Value = $row[$col]
    }

    $hCell = $q1Sheet.Cells.Item($r, 8)
    $hCell.Value = $row[6]
}

# ---------------------------------------------------------------------
# 3. Re-create "总计" right after "2022-Q1" and rebuild its rows, with
#    the new 2022-Q1 totals prepended.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Add($null, $q1Sheet)
$totalSheet.Name = "总计"

$totalHeaders = @("日期", "持有数量(只)", "持有市值(亿元)")
for ($col = 0; $col -lt $totalHeaders.Length; $col++) {
    $cell = $totalSheet.Cells.Item(1, $col + 2)
    $cell.Value = $totalHeaders[$col]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$totalData = @(
    @("2022-Q1", 16, 8.72),
    @("2021-Q4", 7, 3.46),
    @("2021-Q3", 4, 0.68),
    @("2021-Q2", 3, 0.73),
    @("2021-Q1", 4, 3.36),
    @("2020-Q4", 9, 1.99)
)

for ($i = 0; $i -lt $totalData.Length; $i++) {
    $r = $i + 2
    $row = $totalData[$i]

    $aCell = $totalSheet.Cells.Item($r, 1)
    $aCell.Value = $i
    $aCell.Font.Bold = $true
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160
    $aCell.Borders.LineStyle = 1

    $totalSheet.Cells.Item($r, 2).Value = $row[0]
    $totalSheet.Cells.Item($r, 3).Value = $row[1]
    $totalSheet.Cells.Item($r, 4).Value = $row[2]
}
